# Auto-generated edit script applying the scheduled price-refresh update
# to the Leve profit tables across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 5349.737
$ws.Range("I32").Value = 12000
$ws.Range("J32").Value = 4980.278
$ws.Range("K32").Value = 12000
$ws.Range("L32").Value = 4980.278
$ws.Range("M32").Value = -11674
$ws.Range("N32").Value = -5632.278
$ws.Range("H33").Value = 6408
$ws.Range("I33").Value = 7527.909
$ws.Range("J33").Value = 248.5
$ws.Range("K33").Value = 7527.909
$ws.Range("L33").Value = 248.5
$ws.Range("M33").Value = -7298.909
$ws.Range("N33").Value = -706.5
$ws.Range("H53").Value = 4583.48
$ws.Range("I53").Value = 209.75
$ws.Range("K53").Value = 209.75
$ws.Range("M53").Value = 427.25
$ws.Range("H86").Value = 2381.3572
$ws.Range("I86").Value = 2266.2222
$ws.Range("J86").Value = 2588.6
$ws.Range("K86").Value = 2266.2222
$ws.Range("L86").Value = 2588.6
$ws.Range("M86").Value = -1143.2222
$ws.Range("N86").Value = -4834.6
$ws.Range("H89").Value = 2381.3572
$ws.Range("I89").Value = 2266.2222
$ws.Range("J89").Value = 2588.6
$ws.Range("K89").Value = 11331.111
$ws.Range("L89").Value = 12943
$ws.Range("M89").Value = -5715.111000000001
$ws.Range("N89").Value = -24175
$ws.Range("H96").Value = 569.8
$ws.Range("I96").Value = 569.8
$ws.Range("K96").Value = 1709.4
$ws.Range("M96").Value = -336.3999999999999
$ws.Range("H137").Value = 200064.22
$ws.Range("I137").Value = 357155.6
$ws.Range("K137").Value = 1071466.8
$ws.Range("M137").Value = -1068916.8

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 10668.333
$ws.Range("J3").Value = 15000
$ws.Range("L3").Value = 15000
$ws.Range("N3").Value = -15230
$ws.Range("H32").Value = 4317.4546
$ws.Range("I32").Value = 2990.2778
$ws.Range("K32").Value = 2990.2778
$ws.Range("M32").Value = -2703.2778
$ws.Range("H45").Value = 6157449.5
$ws.Range("I45").Value = 10991433
$ws.Range("K45").Value = 10991433
$ws.Range("M45").Value = -10991056
$ws.Range("H74").Value = 54347.562
$ws.Range("I74").Value = 7444.9033
$ws.Range("J74").Value = 236095.38
$ws.Range("K74").Value = 7444.9033
$ws.Range("L74").Value = 236095.38
$ws.Range("M74").Value = -6570.9033
$ws.Range("N74").Value = -237843.38
$ws.Range("H77").Value = 54347.562
$ws.Range("I77").Value = 7444.9033
$ws.Range("J77").Value = 236095.38
$ws.Range("K77").Value = 37224.5165
$ws.Range("L77").Value = 1180476.9
$ws.Range("M77").Value = -32856.5165
$ws.Range("N77").Value = -1189212.9
$ws.Range("H110").Value = 2527611
$ws.Range("I110").Value = 3474761.8
$ws.Range("J110").Value = 1875.3334
$ws.Range("K110").Value = 3474761.8
$ws.Range("L110").Value = 1875.3334
$ws.Range("M110").Value = -3472716.8
$ws.Range("N110").Value = -5965.3334

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 14292713
$ws.Range("I107").Value = 35718284
$ws.Range("J107").Value = 9000
$ws.Range("K107").Value = 35718284
$ws.Range("L107").Value = 9000
$ws.Range("M107").Value = -35716364
$ws.Range("N107").Value = -12840

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1549.4667
$ws.Range("I16").Value = 1396.5
$ws.Range("J16").Value = 2161.3333
$ws.Range("K16").Value = 1396.5
$ws.Range("L16").Value = 2161.3333
$ws.Range("M16").Value = -1109.5
$ws.Range("N16").Value = -2735.3333
$ws.Range("H31").Value = 13746.3955
$ws.Range("I31").Value = 2460.5625
$ws.Range("K31").Value = 2460.5625
$ws.Range("M31").Value = -2165.5625
$ws.Range("H34").Value = 13746.3955
$ws.Range("I34").Value = 2460.5625
$ws.Range("K34").Value = 2460.5625
$ws.Range("M34").Value = -2258.5625
$ws.Range("H58").Value = 3055.8667
$ws.Range("I58").Value = 2089.6667
$ws.Range("K58").Value = 2089.6667
$ws.Range("M58").Value = -1886.6667
$ws.Range("H99").Value = 4217.8335
$ws.Range("I99").Value = 3817
$ws.Range("K99").Value = 3817
$ws.Range("M99").Value = -2319
$ws.Range("H105").Value = 1776.6666
$ws.Range("I105").Value = 1641.4286
$ws.Range("J105").Value = 2250
$ws.Range("K105").Value = 1641.4286
$ws.Range("L105").Value = 2250
$ws.Range("M105").Value = 105.5714
$ws.Range("N105").Value = -5744
$ws.Range("H113").Value = 1549.4667
$ws.Range("I113").Value = 1396.5
$ws.Range("J113").Value = 2161.3333
$ws.Range("K113").Value = 1396.5
$ws.Range("L113").Value = 2161.3333
$ws.Range("M113").Value = 773.5
$ws.Range("N113").Value = -6501.3333
$ws.Range("H126").Value = 4217.8335
$ws.Range("I126").Value = 3817
$ws.Range("K126").Value = 11451
$ws.Range("M126").Value = -8981
$ws.Range("H132").Value = 44597.145
$ws.Range("I132").Value = 2382.6
$ws.Range("K132").Value = 7147.799999999999
$ws.Range("M132").Value = -4617.799999999999
$ws.Range("H134").Value = 2793.28
$ws.Range("I134").Value = 1990.8125
$ws.Range("K134").Value = 5972.4375
$ws.Range("M134").Value = -3437.4375
$ws.Range("H136").Value = 3055.8667
$ws.Range("I136").Value = 2089.6667
$ws.Range("K136").Value = 6269.000100000001
$ws.Range("M136").Value = -3719.000100000001

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 46832.5
$ws.Range("J37").Value = 46832.5
$ws.Range("L37").Value = 140497.5
$ws.Range("N37").Value = -140721.5
$ws.Range("H87").Value = 8832.833000000001
$ws.Range("I87").Value = 8832.833000000001
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 26498.499
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = -25250.499
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 8832.833000000001
$ws.Range("I90").Value = 8832.833000000001
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 79495.497
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = -73255.497
$ws.Range("N90").ClearContents()
$ws.Range("H113").Value = 3294.25
$ws.Range("I113").Value = 6979
$ws.Range("J113").Value = 2324.5789
$ws.Range("K113").Value = 20937
$ws.Range("L113").Value = 6973.736699999999
$ws.Range("M113").Value = -18767
$ws.Range("N113").Value = -11313.7367
$ws.Range("H140").Value = 2929.6667
$ws.Range("I140").Value = 2795.875
$ws.Range("J140").Value = 4000
$ws.Range("K140").Value = 8387.625
$ws.Range("L140").Value = 12000
$ws.Range("M140").Value = -3207.625
$ws.Range("N140").Value = -22360

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6748825
$ws.Range("I102").Value = 9261457
$ws.Range("J102").Value = 2979876.8
$ws.Range("K102").Value = 9261457
$ws.Range("L102").Value = 2979876.8
$ws.Range("M102").Value = -9259835
$ws.Range("N102").Value = -2983120.8
$ws.Range("H126").Value = 8588937
$ws.Range("I126").Value = 4548483
$ws.Range("K126").Value = 13645449
$ws.Range("M126").Value = -13642979
$ws.Range("H132").Value = 3772.0312
$ws.Range("I132").Value = 3130.1904
$ws.Range("J132").Value = 4997.364
$ws.Range("K132").Value = 9390.5712
$ws.Range("L132").Value = 14992.092
$ws.Range("M132").Value = -6860.5712
$ws.Range("N132").Value = -20052.092

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 82426.09
$ws.Range("I22").Value = 444944
$ws.Range("K22").Value = 444944
$ws.Range("M22").Value = -444649
$ws.Range("H27").Value = 82426.09
$ws.Range("I27").Value = 444944
$ws.Range("K27").Value = 444944
$ws.Range("M27").Value = -444837
$ws.Range("H61").Value = 4632145
$ws.Range("I61").Value = 5293755
$ws.Range("K61").Value = 5293755
$ws.Range("M61").Value = -5293553
$ws.Range("H113").Value = 4632145
$ws.Range("I113").Value = 5293755
$ws.Range("K113").Value = 5293755
$ws.Range("M113").Value = -5291585
$ws.Range("H122").Value = 5350.522
$ws.Range("I122").Value = 3951.6
$ws.Range("K122").Value = 11854.8
$ws.Range("M122").Value = -9404.799999999999

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2152.8572
$ws.Range("I122").Value = 1415
$ws.Range("K122").Value = 4245
$ws.Range("M122").Value = -1795
$ws.Range("H126").Value = 4181.727
$ws.Range("I126").Value = 3999.889
$ws.Range("K126").Value = 11999.667
$ws.Range("M126").Value = -9529.667000000001

